# Append 36 new translation rows (my-properties feature) to Sheet1
# Columns: A=key, B=en, C=hi, D=te
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 36,4
$arr[0,0] = 'my_properties_title'
$arr[0,1] = 'My Properties'
$arr[0,2] = 'मेरी संपत्तियां'
$arr[0,3] = 'నా ప్రాపర్టీలు'
$arr[1,0] = 'my_properties_search_placeholder'
$arr[1,1] = 'Search properties...'
$arr[1,2] = 'संपत्तियां खोजें...'
$arr[1,3] = 'ప్రాపర్టీలు వెతకండి...'
$arr[2,0] = 'tab_all'
$arr[2,1] = 'All'
$arr[2,2] = 'सभी'
$arr[2,3] = 'అన్నీ'
$arr[3,0] = 'tab_active'
$arr[3,1] = 'Active'
$arr[3,2] = 'सक्रिय'
$arr[3,3] = 'యాక్టివ్'
$arr[4,0] = 'tab_pending'
$arr[4,1] = 'Pending'
$arr[4,2] = 'लंबित'
$arr[4,3] = 'పెండింగ్'
$arr[5,0] = 'tab_sold'
$arr[5,1] = 'Sold'
$arr[5,2] = 'बेचा गया'
$arr[5,3] = 'అమ్మబడింది'
$arr[6,0] = 'status_active'
$arr[6,1] = 'Active'
$arr[6,2] = 'सक्रिय'
$arr[6,3] = 'యాక్టివ్'
$arr[7,0] = 'status_pending_review'
$arr[7,1] = 'Pending Review'
$arr[7,2] = 'समीक्षा लंबित'
$arr[7,3] = 'సమీక్ష పెండింగ్'
$arr[8,0] = 'status_rejected'
$arr[8,1] = 'Rejected'
$arr[8,2] = 'अस्वीकृत'
$arr[8,3] = 'తిరస్కరించబడింది'
$arr[9,0] = 'status_sold_out'
$arr[9,1] = 'Sold Out'
$arr[9,2] = 'बिक गया'
$arr[9,3] = 'అమ్మబడింది'
$arr[10,0] = 'edit_to_resubmit'
$arr[10,1] = 'Edit to resubmit'
$arr[10,2] = 'पुनः सबमिट करने के लिए संपादित करें'
$arr[10,3] = 'మళ్ళీ సమర్పించడానికి సవరించండి'
$arr[11,0] = 'no_reviews_yet'
$arr[11,1] = 'No reviews yet'
$arr[11,2] = 'अभी कोई समीक्षा नहीं'
$arr[11,3] = 'ఇంకా సమీక్షలు లేవు'
$arr[12,0] = 'verified'
$arr[12,1] = 'Verified'
$arr[12,2] = 'सत्यापित'
$arr[12,3] = 'ధృవీకరించబడింది'
$arr[13,0] = 'loading_properties'
$arr[13,1] = 'Loading properties...'
$arr[13,2] = 'संपत्तियां लोड हो रही हैं...'
$arr[13,3] = 'ప్రాపర్టీలు లోడ్ అవుతున్నాయి...'
$arr[14,0] = 'no_properties_found'
$arr[14,1] = 'No properties found'
$arr[14,2] = 'कोई संपत्ति नहीं मिली'
$arr[14,3] = 'ప్రాపర్టీలు కనుగొనబడలేదు'
$arr[15,0] = 'no_properties_added'
$arr[15,1] = 'You haven''t added any properties yet'
$arr[15,2] = 'आपने अभी तक कोई संपत्ति नहीं जोड़ी'
$arr[15,3] = 'మీరు ఇంకా ఏ ప్రాపర్టీలు జోడించలేదు'
$arr[16,0] = 'no_properties_in_tab'
$arr[16,1] = 'No {tab} properties'
$arr[16,2] = 'कोई {tab} संपत्ति नहीं'
$arr[16,3] = '{tab} ప్రాపర్టీలు లేవు'
$arr[17,0] = 'search_no_match'
$arr[17,1] = 'No properties match "{query}"'
$arr[17,2] = '{query} से मेल खाती कोई संपत्ति नहीं'
$arr[17,3] = '{query} కు సరిపోలే ప్రాపర్టీలు లేవు'
$arr[18,0] = 'untitled_property'
$arr[18,1] = 'Untitled Property'
$arr[18,2] = 'बिना शीर्षक की संपत्ति'
$arr[18,3] = 'శీర్షిక లేని ప్రాపర్టీ'
$arr[19,0] = 'location_not_specified'
$arr[19,1] = 'Location not specified'
$arr[19,2] = 'स्थान निर्दिष्ट नहीं'
$arr[19,3] = 'స్థానం పేర్కొనబడలేదు'
$arr[20,0] = 'btn_view'
$arr[20,1] = 'View'
$arr[20,2] = 'देखें'
$arr[20,3] = 'చూడండి'
$arr[21,0] = 'options_title'
$arr[21,1] = 'Property Options'
$arr[21,2] = 'संपत्ति विकल्प'
$arr[21,3] = 'ప్రాపర్టీ ఎంపికలు'
$arr[22,0] = 'option_edit'
$arr[22,1] = 'Edit Property'
$arr[22,2] = 'संपत्ति संपादित करें'
$arr[22,3] = 'ప్రాపర్టీ సవరించండి'
$arr[23,0] = 'option_cannot_edit_sold'
$arr[23,1] = 'Cannot Edit (Sold)'
$arr[23,2] = 'संपादित नहीं कर सकते (बेचा गया)'
$arr[23,3] = 'సవరించలేరు (అమ్మబడింది)'
$arr[24,0] = 'option_resubmit_note'
$arr[24,1] = 'Property will be resubmitted for review'
$arr[24,2] = 'संपत्ति समीक्षा के लिए पुनः सबमिट की जाएगी'
$arr[24,3] = 'ప్రాపర్టీ సమీక్షకు మళ్ళీ సమర్పించబడుతుంది'
$arr[25,0] = 'option_mark_sold'
$arr[25,1] = 'Mark as Sold'
$arr[25,2] = 'बेचा गया चिह्नित करें'
$arr[25,3] = 'అమ్మినట్లు గుర్తించండి'
$arr[26,0] = 'option_already_sold'
$arr[26,1] = 'Already Sold'
$arr[26,2] = 'पहले से बेचा गया'
$arr[26,3] = 'ఇప్పటికే అమ్మబడింది'
$arr[27,0] = 'btn_cancel'
$arr[27,1] = 'Cancel'
$arr[27,2] = 'रद्द करें'
$arr[27,3] = 'రద్దు చేయి'
$arr[28,0] = 'mark_sold_title'
$arr[28,1] = 'Mark as Sold'
$arr[28,2] = 'बेचा गया चिह्नित करें'
$arr[28,3] = 'అమ్మినట్లు గుర్తించండి'
$arr[29,0] = 'mark_sold_message'
$arr[29,1] = 'Are you sure you want to mark this property as sold? This action cannot be undone by you.'
$arr[29,2] = 'क्या आप वाकई इस संपत्ति को बेचा गया चिह्नित करना चाहते हैं? यह क्रिया आप पूर्ववत नहीं कर सकते।'
$arr[29,3] = 'మీరు ఈ ప్రాపర్టీని అమ్మినట్లు గుర్తించాలనుకుంటున్నారా? ఈ చర్యను మీరు రద్దు చేయలేరు.'
$arr[30,0] = 'mark_sold_confirm'
$arr[30,1] = 'Mark as Sold'
$arr[30,2] = 'बेचा गया चिह्नित करें'
$arr[30,3] = 'అమ్మినట్లు గుర్తించండి'
$arr[31,0] = 'success_marked_sold'
$arr[31,1] = 'Property marked as sold'
$arr[31,2] = 'संपत्ति बेचा गया चिह्नित की गई'
$arr[31,3] = 'ప్రాపర్టీ అమ్మినట్లు గుర్తించబడింది'
$arr[32,0] = 'error_mark_sold'
$arr[32,1] = 'Failed to mark property as sold'
$arr[32,2] = 'संपत्ति को बेचा गया चिह्नित करना विफल रहा'
$arr[32,3] = 'ప్రాపర్టీని అమ్మినట్లు గుర్తించడం విఫలమైంది'
$arr[33,0] = 'voice_search_not_supported'
$arr[33,1] = 'Voice search is not supported in this browser'
$arr[33,2] = 'इस ब्राउज़र में वॉइस सर्च समर्थित नहीं है'
$arr[33,3] = 'ఈ బ్రౌజర్‌లో వాయిస్ సెర్చ్ మద్దతు లేదు'
$arr[34,0] = 'review_singular'
$arr[34,1] = 'review'
$arr[34,2] = 'समीक्षा'
$arr[34,3] = 'సమీక్ష'
$arr[35,0] = 'review_plural'
$arr[35,1] = 'reviews'
$arr[35,2] = 'समीक्षाएं'
$arr[35,3] = 'సమీక్షలు'

$startRow = 1715
$endRow = 1750
$rng = $ws.Range("A" + $startRow + ":D" + $endRow)
$rng.Value2 = $arr

# Match the author-recorded view state (scroll position + active cell)
$ws.Range("A1715").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1691
$excel.ActiveWindow.ScrollColumn = 1